$wb = $excel.ActiveWorkbook

# Rename "Mapping Specification" sheet to "Mappings"
$wsMappings = $wb.Worksheets.Item("Mapping Specification")
$wsMappings.Name = "Mappings"

# Commit message: content on "Templates" sheet was used to create a HumanMessage;
# the author was reviewing the Templates sheet (selection moved there), then
# returned focus to the Mappings sheet, which is now the active/selected tab.
$wsTemplates = $wb.Worksheets.Item("Templates")
$wsTemplates.Activate()
$wsTemplates.Range("B2").Select()

$wsMappings.Activate()
$wsMappings.Range("A43").Select()
